$p = $ppt.ActivePresentation

# Slide 1 ("API Review: Deprecation..." slide), shape 2 ("Text Placeholder 8"):
# drop the "Ruslan Israfilov" paragraph, keep only "ww12'21"
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.Delete()
